$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "ADD CRUD in the first Window" - append a new record row (row 4) to the
# table of judges/lawyers, mirroring the existing rows.

# RUC / ID Number for this record are numeric-looking strings that must be
# stored as text (like row 3's C3/D3), so force a text number format before
# assigning the values.
$ws.Range("C4:D4").NumberFormat = "@"

$ws.Range("A4").Value = "Jesus"
$ws.Range("B4").Value = "Alexander Benitez"
$ws.Range("C4").Value = "1458789"
$ws.Range("D4").Value = "1728224557"
$ws.Range("E4").Value = "Dr. Atiencia Atiencia Atiencia Atiencia"
